$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 8 (A8..I8) down to row 9 first so the new row
# matches the existing style (date format in col A, boolean style in col G, etc.)
$ws.Range("A8:I8").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)

$ws.Range("A9").Value = 42654.745717592596
$ws.Range("B9").Value = $false
$ws.Range("C9").Value = 10121.91
$ws.Range("D9").Value = 10126.469999999999
$ws.Range("E9").Value = 75.5
$ws.Range("F9").Value = 75.569999999999993
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = 0.09
$ws.Range("I9").Value = $false

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null
$ws.Columns.Item(7).AutoFit() | Out-Null
$ws.Columns.Item(8).AutoFit() | Out-Null
$ws.Columns.Item(9).AutoFit() | Out-Null
